$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10599.8
$ws.Range("H106").Value = 10000
$ws.Range("J106").Value = 10000
$ws.Range("L106").Value = 10000
$ws.Range("N106").Value = -11262
$ws.Range("H132").Value = 907.61536
$ws.Range("I132").Value = 907.61536
$ws.Range("K132").Value = 2722.84608
$ws.Range("M132").Value = -192.8460800000003
$ws.Range("H136").Value = 75000
$ws.Range("J136").Value = 75000
$ws.Range("L136").Value = 75000
$ws.Range("N136").Value = -85200
$ws.Range("H137").Value = 2490.9697
$ws.Range("J137").Value = 4079.4614
$ws.Range("L137").Value = 12238.3842
$ws.Range("N137").Value = -17338.3842
$ws.Range("H138").Value = 4332.5186
$ws.Range("I138").Value = 4016.818
$ws.Range("J138").Value = 4549.5625
$ws.Range("K138").Value = 12050.454
$ws.Range("L138").Value = 13648.6875
$ws.Range("M138").Value = -6910.454000000002
$ws.Range("N138").Value = -23928.6875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6198.096
$ws.Range("I32").Value = 5930.575
$ws.Range("K32").Value = 5930.575
$ws.Range("M32").Value = -5643.575
$ws.Range("H61").Value = 1479.2413
$ws.Range("I61").Value = 1024.8636
$ws.Range("K61").Value = 1024.8636
$ws.Range("M61").Value = -812.8635999999999
$ws.Range("H74").Value = 1973.1305
$ws.Range("I74").Value = 1298.9412
$ws.Range("J74").Value = 3883.3333
$ws.Range("K74").Value = 1298.9412
$ws.Range("L74").Value = 3883.3333
$ws.Range("M74").Value = -424.9412
$ws.Range("N74").Value = -5631.3333
$ws.Range("H77").Value = 1973.1305
$ws.Range("I77").Value = 1298.9412
$ws.Range("J77").Value = 3883.3333
$ws.Range("K77").Value = 6494.706
$ws.Range("L77").Value = 19416.6665
$ws.Range("M77").Value = -2126.706
$ws.Range("N77").Value = -28152.6665
$ws.Range("H119").Value = 90000
$ws.Range("J119").Value = 90000
$ws.Range("L119").Value = 90000
$ws.Range("N119").Value = -99676
$ws.Range("H132").Value = 2672.6553
$ws.Range("I132").Value = 1925.6
$ws.Range("K132").Value = 5776.799999999999
$ws.Range("M132").Value = -3246.799999999999
$ws.Range("H136").Value = 1479.2413
$ws.Range("I136").Value = 1024.8636
$ws.Range("K136").Value = 3074.5908
$ws.Range("M136").Value = -524.5907999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4196.3335
$ws.Range("I86").Value = 3320.8
$ws.Range("K86").Value = 3320.8
$ws.Range("M86").Value = -2197.8
$ws.Range("H89").Value = 4196.3335
$ws.Range("I89").Value = 3320.8
$ws.Range("K89").Value = 16604
$ws.Range("M89").Value = -10988
$ws.Range("H134").Value = 2237.9722
$ws.Range("I134").Value = 1816
$ws.Range("J134").Value = 4347.8335
$ws.Range("K134").Value = 5448
$ws.Range("L134").Value = 13043.5005
$ws.Range("M134").Value = -2913
$ws.Range("N134").Value = -18113.5005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4602.1665
$ws.Range("I25").Value = 1520
$ws.Range("K25").Value = 1520
$ws.Range("M25").Value = -1346
$ws.Range("H31").Value = 2417.7058
$ws.Range("I31").Value = 2421.75
$ws.Range("K31").Value = 2421.75
$ws.Range("M31").Value = -2126.75
$ws.Range("H34").Value = 2417.7058
$ws.Range("I34").Value = 2421.75
$ws.Range("K34").Value = 2421.75
$ws.Range("M34").Value = -2219.75
$ws.Range("H58").Value = 2918
$ws.Range("I58").Value = 2848.7273
$ws.Range("K58").Value = 2848.7273
$ws.Range("M58").Value = -2645.7273
$ws.Range("H122").Value = 1633
$ws.Range("I122").Value = 1727.5
$ws.Range("K122").Value = 5182.5
$ws.Range("M122").Value = -2732.5
$ws.Range("H132").Value = 4244.1904
$ws.Range("I132").Value = 4087.7334
$ws.Range("K132").Value = 12263.2002
$ws.Range("M132").Value = -9733.200199999999
$ws.Range("H134").Value = 2556.5
$ws.Range("I134").Value = 2368.6155
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 7105.8465
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -4570.8465
$ws.Range("N134").Value = -20067
$ws.Range("H136").Value = 2918
$ws.Range("I136").Value = 2848.7273
$ws.Range("K136").Value = 8546.1819
$ws.Range("M136").Value = -5996.1819
$ws.Range("H141").Value = 52057.25
$ws.Range("J141").Value = 52057.25
$ws.Range("L141").Value = 52057.25
$ws.Range("N141").Value = -62417.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 101699.6
$ws.Range("J122").Value = 101699.6
$ws.Range("L122").Value = 915296.4
$ws.Range("N122").Value = -920196.4
$ws.Range("H131").Value = 1526.8572
$ws.Range("I131").Value = 1052.8334
$ws.Range("J131").Value = 1882.375
$ws.Range("K131").Value = 3158.5002
$ws.Range("L131").Value = 5647.125
$ws.Range("M131").Value = 1881.4998
$ws.Range("N131").Value = -15727.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 20502
$ws.Range("J18").Value = 20502
$ws.Range("L18").Value = 20502
$ws.Range("N18").Value = -21088
$ws.Range("H46").Value = 9608.75
$ws.Range("J46").Value = 9608.75
$ws.Range("L46").Value = 9608.75
$ws.Range("N46").Value = -9920.75
$ws.Range("H57").Value = 19633.334
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 19633.334
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 19633.334
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21273.334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1999
$ws.Range("I7").Value = 1999
$ws.Range("J7").Value = 1999
$ws.Range("K7").Value = 1999
$ws.Range("L7").Value = 1999
$ws.Range("M7").Value = -1887
$ws.Range("N7").Value = -2223
$ws.Range("H55").Value = 397.42856
$ws.Range("J55").Value = 379.5
$ws.Range("L55").Value = 379.5
$ws.Range("N55").Value = -725.5
$ws.Range("H93").Value = 1933.5
$ws.Range("I93").Value = 1933.5
$ws.Range("K93").Value = 1933.5
$ws.Range("M93").Value = -685.5
$ws.Range("H126").Value = 1999
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 5997
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -10937
$ws.Range("H132").Value = 5229.4
$ws.Range("I132").Value = 5037
$ws.Range("K132").Value = 15111
$ws.Range("M132").Value = -12581
$ws.Range("H136").Value = 3599.8064
$ws.Range("I136").Value = 3903.0417
$ws.Range("J136").Value = 2560.1428
$ws.Range("K136").Value = 11709.1251
$ws.Range("L136").Value = 7680.428400000001
$ws.Range("M136").Value = -9159.125100000001
$ws.Range("N136").Value = -12780.4284
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 3743.8667
$ws.Range("I81").Value = 2654
$ws.Range("J81").Value = 19002
$ws.Range("K81").Value = 5308
$ws.Range("L81").Value = 38004
$ws.Range("M81").Value = -4247
$ws.Range("N81").Value = -40126
$ws.Range("H84").Value = 3743.8667
$ws.Range("I84").Value = 2654
$ws.Range("J84").Value = 19002
$ws.Range("K84").Value = 26540
$ws.Range("L84").Value = 190020
$ws.Range("M84").Value = -21236
$ws.Range("N84").Value = -200628
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H107").Value = 542.2857
$ws.Range("I107").Value = 499.33334
$ws.Range("K107").Value = 1498.00002
$ws.Range("M107").Value = 421.9999800000001
$ws.Range("H126").Value = 2184.8667
$ws.Range("I126").Value = 2182.6155
$ws.Range("J126").Value = 2199.5
$ws.Range("K126").Value = 6547.8465
$ws.Range("L126").Value = 6598.5
$ws.Range("M126").Value = -4077.8465
$ws.Range("N126").Value = -11538.5
$ws.Range("H132").Value = 4000.5
$ws.Range("I132").Value = 3002
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9006
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6476
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 3408.0715
$ws.Range("I136").Value = 3207.3684
$ws.Range("K136").Value = 9622.1052
$ws.Range("M136").Value = -7072.1052
